# feat: add 2022-Q1 data
#
# Inserts a new "2022-Q1" sheet (fund-holdings detail, same layout as the
# other quarterly sheets) right before the "总计" (totals) summary sheet,
# and prepends a matching "2022-Q1" row to the "总计" sheet.
#
# NOTE: worksheet/range/cell variables in this host resolve by POSITION
# at the time they are dereferenced, not by stable object identity, so we
# deliberately re-fetch sheets/ranges by name/index right before each use
# instead of reusing handles captured before a structural change (sheet
# insert, row insert, rename, ...).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Create the new "2022-Q1" sheet positioned right before "总计"
# ---------------------------------------------------------------------
$lastIndex = $wb.Worksheets.Count
$newSheet  = $wb.Worksheets.Add($wb.Worksheets.Item($lastIndex))   # insert before "总计" (last sheet)
$newSheet.Name = "2022-Q1"

# Carry over the header-row / index-column formatting (style index,
# borders, bold, alignment, ...) from an existing quarterly sheet so the
# new sheet matches the workbook's established look.
$wb.Worksheets.Item("2021-Q4").Range("B1:H1").Copy()
$wb.Worksheets.Item("2022-Q1").Range("B1:H1").PasteSpecial(-4122)   # xlPasteFormats

$wb.Worksheets.Item("2021-Q4").Range("A2:A4").Copy()
$wb.Worksheets.Item("2022-Q1").Range("A2:A4").PasteSpecial(-4122)   # xlPasteFormats

$excel.CutCopyMode = 0

$q1 = $wb.Worksheets.Item("2022-Q1")

# Header row
$q1.Cells.Item(1,2).Value = "基金代码"
$q1.Cells.Item(1,3).Value = "基金名称"
$q1.Cells.Item(1,4).Value = "基金规模"
$q1.Cells.Item(1,5).Value = "股票总仓位"
$q1.Cells.Item(1,6).Value = "仓位占比"
$q1.Cells.Item(1,7).Value = "持有市值(亿元)"
$q1.Cells.Item(1,8).Value = "仓位排名"

# Index column (A2:A4) -- plain integers, 0-based row counter
$q1.Cells.Item(2,1).Value = 0
$q1.Cells.Item(3,1).Value = 1
$q1.Cells.Item(4,1).Value = 2

# Columns B (fund code) and D/E/F/G hold numeric-looking values that are
# stored as *text* in this workbook (leading zeros / fixed decimals must be
# preserved), so force a text number-format on them before writing values.
# (Two statements on purpose -- a single comma-separated multi-area range
# only picks up the format on its first area in this host.)
$q1.Range("B2:B4").NumberFormat = "@"
$q1.Range("D2:G4").NumberFormat = "@"

# Row 2 - 011383 / 富安达医药创新混合
$q1.Cells.Item(2,2).Value = "011383"
$q1.Cells.Item(2,3).Value = "富安达医药创新混合"
$q1.Cells.Item(2,4).Value = "1.68"
$q1.Cells.Item(2,5).Value = "83.50"
$q1.Cells.Item(2,6).Value = "4.62"
$q1.Cells.Item(2,7).Value = "0.0776"
$q1.Cells.Item(2,8).Value = 2

# Row 3 - 001861 / 富安达健康人生灵活配置混合
$q1.Cells.Item(3,2).Value = "001861"
$q1.Cells.Item(3,3).Value = "富安达健康人生灵活配置混合"
$q1.Cells.Item(3,4).Value = "0.61"
$q1.Cells.Item(3,5).Value = "82.18"
$q1.Cells.Item(3,6).Value = "3.82"
$q1.Cells.Item(3,7).Value = "0.0233"
$q1.Cells.Item(3,8).Value = 6

# Row 4 - 001972 / 前海开源沪港深智慧生活优选灵活配置混合
$q1.Cells.Item(4,2).Value = "001972"
$q1.Cells.Item(4,3).Value = "前海开源沪港深智慧生活优选灵活配置混合"
$q1.Cells.Item(4,4).Value = "0.48"
$q1.Cells.Item(4,5).Value = "77.26"
$q1.Cells.Item(4,6).Value = "3.71"
$q1.Cells.Item(4,7).Value = "0.0178"
$q1.Cells.Item(4,8).Value = 9

# The "@" text format above was only needed to stop Excel from silently
# re-parsing these strings as numbers (losing the leading zero / trailing
# zero digits) while the values were written; drop it again afterwards so
# the cells end up unstyled, same as the other quarterly sheets.
$q1.Range("B2:B4").ClearFormats()
$q1.Range("D2:G4").ClearFormats()

# ---------------------------------------------------------------------
# 2. Update "总计" - insert a new top data row for 2022-Q1 and shift the
#    existing rows (and their index column) down by one.
# ---------------------------------------------------------------------
$wb.Worksheets.Item("总计").Rows.Item(2).Insert()

# Row-insert clones the formatting of the row above (the bold/bordered
# header row) onto the whole new blank row; the data rows in this sheet
# only carry an explicit style on column A, so reset B:D and then copy the
# correct "index column" style (from the row that used to be row 2, now
# row 3) onto the new A2.
$wb.Worksheets.Item("总计").Range("A2:D2").ClearFormats()
$wb.Worksheets.Item("总计").Range("A3").Copy()
$wb.Worksheets.Item("总计").Range("A2").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

$total = $wb.Worksheets.Item("总计")
$total.Cells.Item(2,1).Value = 0
$total.Cells.Item(2,2).Value = "2022-Q1"
$total.Cells.Item(2,3).Value = 3
$total.Cells.Item(2,4).Value = 0.12

# Re-number the index column (A) for the rows that shifted down
$total.Cells.Item(3,1).Value = 1
$total.Cells.Item(4,1).Value = 2
$total.Cells.Item(5,1).Value = 3
$total.Cells.Item(6,1).Value = 4
$total.Cells.Item(7,1).Value = 5
